# Add Japanese translations (column G) aligned with the English terms in column B
# for rows 2-82, matching the order/content added to the workbook upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$translations = @(
    "寝汗",
    "労作後症状の増悪",
    "壊疽",
    "精神的健康の低下",
    "身体機能の低下",
    "てんかん性突然死",
    "心臓突然死",
    "無力症",
    "冷え性",
    "硬直",
    "肝不全",
    "再発性単発性肝炎",
    "発作性直腸痛",
    "肩甲骨内痛",
    "腰痛",
    "坐骨神経痛",
    "突破痛",
    "発作性腹痛",
    "脇腹痛",
    "腱痛",
    "股関節痛",
    "腹部硬直",
    "腹部けいれん",
    "消化不良",
    "腹部疝痛",
    "腹部固縮",
    "急性腹症",
    "骨痛",
    "陰嚢痛",
    "舌痛症",
    "首の痛み",
    "眼球運動誘発性疼痛",
    "耳の痛み",
    "慢性疼痛",
    "胸膜炎性胸痛",
    "前胸部痛",
    "下顎痛",
    "骨盤痛",
    "肩の痛み",
    "外陰部痛",
    "膀胱痛",
    "赤芽球痛",
    "精巣痛",
    "関節痛/関節炎",
    "アロディニア",
    "乳房痛",
    "ウィルソン徴候",
    "足の痛み",
    "足指痛",
    "足首の痛み",
    "手首の痛み",
    "肘の痛み",
    "指の痛み",
    "顎の痛み",
    "肋軟骨痛",
    "後頭神経痛",
    "三叉神経痛",
    "心窩部痛",
    "咽頭痛",
    "鼠径部痛",
    "腓腹筋痛",
    "運動誘発性筋肉痛",
    "ほてり",
    "虚弱",
    "倦怠感",
    "移乗能力の低下",
    "入浴能力の低下",
    "摂食能力の低下",
    "尿失禁",
    "ストレス性尿失禁",
    "着替え能力の障害",
    "排泄能力の障害",
    "膣内の魚臭",
    "かび臭さ",
    "職業障害",
    "体の痛み",
    "胸の張り",
    "発作性のくしゃみ",
    "認知疲労",
    "慢性疲労",
    "周期性疲労"
)

for ($i = 0; $i -lt $translations.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $translations[$i]
}

# Mirror the header-row selection state left behind by the edit (whole first row selected)
$ws.Rows.Item(1).Select() | Out-Null

